$wb = $excel.ActiveWorkbook

$changes = @{
    "F2"  = 1554
    "F3"  = 8865
    "F7"  = 322
    "F9"  = 33
    "F10" = 44
    "F11" = 3746
    "F14" = 97
    "F15" = 3893
    "F18" = 325
    "F20" = 2534
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $changes.Keys) {
        $ws.Range($addr).Value = $changes[$addr]
    }
}
